$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 269; existing rows 269:341 shift down to 270:342.
$ws.Rows(269).Insert()

# Populate the newly inserted row 269 with the new record
# (Apio, Macroferia Regional de Talca, Maule - week of 2023-07-28).
$ws.Cells.Item(269, 1).Value = 5
$ws.Cells.Item(269, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(269, 3).Value = "Maule"
$ws.Cells.Item(269, 4).Value = 45135
$ws.Cells.Item(269, 5).Value = 7
$ws.Cells.Item(269, 6).Value = 100112017
$ws.Cells.Item(269, 7).Value = "Apio"
$ws.Cells.Item(269, 8).Value = "Americana (o)"
$ws.Cells.Item(269, 9).Value = "Primera"
$ws.Cells.Item(269, 10).Value = 700
$ws.Cells.Item(269, 11).Value = 5500
$ws.Cells.Item(269, 12).Value = 5500
$ws.Cells.Item(269, 13).Value = 5500
$ws.Cells.Item(269, 14).Value = "`$/docena de matas"
$ws.Cells.Item(269, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(269, 16).Value = 917
$ws.Cells.Item(269, 17).Value = 6
$ws.Cells.Item(269, 18).Value = "Hortaliza"
